# conductor_diagnostic.xlsx - "update input file of TDD case1 to last version"
#
# The two worksheets were renamed to match the new naming convention used
# across the rest of the TDD_examples workbooks:
#   "Space" -> "Spatial_distribution"
#   "Time"  -> "Time_evolution"
#
# Excel automatically repoints every in-workbook formula that referenced the
# old sheet names (e.g. "=Space!A1" on the Time sheet) to the new names, so
# no formula text needs to be touched by hand.
#
# The commit also leaves the two sheets with a fresh selection: the
# (now) "Spatial_distribution" sheet ends up with A3:D3 selected/active and
# stays the active tab, while "Time_evolution" is left with A3 selected.

$wb = $excel.ActiveWorkbook

$wsSpace = $wb.Worksheets.Item("Space")
$wsTime  = $wb.Worksheets.Item("Time")

$wsSpace.Name = "Spatial_distribution"
$wsTime.Name  = "Time_evolution"

# Set the selection/active-cell state on "Time_evolution" first ...
[void]$wsTime.Activate()
[void]$wsTime.Range("A3").Select()

# ... then finish on "Spatial_distribution" so it is the tab left active,
# matching tabSelected="1" on that sheet in the saved workbook.
[void]$wsSpace.Activate()
[void]$wsSpace.Range("A3:D3").Select()
